$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.0609130859375
$ws.Range("B1").Value = 1.769896388053894
$ws.Range("C1").Value = 5.206596374511719
$ws.Range("D1").Value = 0.872567892074585
$ws.Range("E1").Value = 0.4049972891807556
